# Add a new "Role-Based UI Visibility Update" status section at the end of the
# document, mirroring the existing section pattern (blank line, "---",
# title, "Updated: <date>", blank line, header row, then the data row).
$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">---</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Role-Based UI Visibility Update</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Updated: 2026-02-18</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Module Name</w:t><w:tab/><w:t xml:space="preserve">Developed</w:t><w:tab/><w:t xml:space="preserve">Partial Developed</w:t><w:tab/><w:t xml:space="preserve">Need To Develop</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Platform vs Store Feature Visibility</w:t><w:tab/><w:t xml:space="preserve">Added backend access summary endpoint (`GET /api/auth/access`) and dynamic sidebar visibility by platform/store scope so platform owner/staff and store users no longer see identical menu options</w:t><w:tab/><w:t xml:space="preserve">Page-level button/action guards still partly API-enforced (403 fallback) and not fully hidden in all screens</w:t><w:tab/><w:t xml:space="preserve">Complete route/action-level guard matrix and per-component permission gating for every admin view</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$insertPoint.InsertXML($xmlFragment)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
Write-Output ("Last paragraph text=" + $d.Paragraphs.Last.Range.Text)
